$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the placeholder values to use templating braces.
$ws.Range("B1").Value = "{firstname}"
$ws.Range("G1").Value = "{{expenses}}"
$ws.Range("B3").Value = "{{address}}"
$ws.Range("B5").Value = "{{hobbies}}"

# Move the active selection to I7.
$ws.Range("I7").Select()
